# Rename the sole worksheet from "Sheet1" to "Data", as happened when the
# author re-uploaded/re-saved the workbook to GitHub.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "Data"
